$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.237.76'
$ws.Range("E2").Value = '  +1.56%  '

$ws.Range("D3").Value = '2.006.47'
$ws.Range("E3").Value = '  +2.20%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.55'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.07%  '

$ws.Range("E6").Value = '  +1.41%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '60.12'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.72%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("E9").Value = '  +2.19%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0810'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.95%  '

$ws.Range("E11").Value = '  +0.74%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.08'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.99%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.48'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.31%  '

$ws.Range("D14").Value = '2.299.55'
$ws.Range("E14").Value = '  +2.11%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.845'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.32%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.45'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.76%  '

$ws.Range("D17").Value = '2.006.09'
$ws.Range("E17").Value = '  +2.24%  '

$ws.Range("D18").Value = '37.170.08'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '70.26'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.53%  '

$ws.Range("E20").Value = '  +1.32%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.21'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.25%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '230.47'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.01%  '

$ws.Range("E23").Value = '  +0.01%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.48'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.31%  '

$ws.Range("E25").Value = '  +0.55%  '

$ws.Range("E26").Value = '  +2.14%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '164.65'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.33%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.138'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.78%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.68'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.82%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.39'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +14.01%  '

$ws.Range("E31").Value = '  +1.24%  '

$ws.Range("E32").Value = '  +0.67%  '

$ws.Range("E33").Value = '  +5.82%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.49'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.11%  '

$ws.Range("E35").Value = '  +5.57%  '

$ws.Range("E36").Value = '  -0.14%  '

$ws.Range("E37").Value = '  +2.28%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.36'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.77%  '

$ws.Range("E39").Value = '  -3.99%  '

$ws.Range("E40").Value = '  +0.21%  '

$ws.Range("E41").Value = '  +0.82%  '

$ws.Range("E42").Value = '  +0.56%  '

$ws.Range("E43").Value = '  +1.40%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.65'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.93%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '91.98'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.92%  '

$ws.Range("D46").Value = '1.370.74'
$ws.Range("E46").Value = '  +0.01%  '

$ws.Range("E47").Value = '  +1.40%  '

$ws.Range("E48").Value = '  +2.57%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.10'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +13.64%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '46.99'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.56%  '

$ws.Range("E51").Value = '  -0.29%  '
